# Update example-import-file.xlsx to the new "import order" sample data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = "Member"
$ws.Range("F2").Value = "Sukabumi"
$ws.Range("G2").Value = "Gunung Puyuh"
$ws.Range("H2").Value = "Alamat 1"
$ws.Range("I2").Value = 43123
$ws.Range("J2").Value = "0000000001"
$ws.Range("L2").Value = "pcs"
$ws.Range("M2").Value = "sepatu"
$ws.Range("R2").Value = "IDR"
$ws.Range("T2").Value = "test note"
$ws.Range("U2").Value = "test reference"

# Row 3
$ws.Range("A3").Value = "1234567891"
$ws.Range("B3").Value = "Member"
$ws.Range("F3").Value = "Sukabumi"
$ws.Range("G3").Value = "Gunung Puyuh"
$ws.Range("H3").Value = "Alamat 1"
$ws.Range("I3").Value = 43123
$ws.Range("J3").Value = "0000000001"
$ws.Range("L3").Value = "pcs"
$ws.Range("M3").Value = "adidas"
$ws.Range("R3").Value = "IDR"
$ws.Range("T3").Value = "test note 2"
$ws.Range("U3").Value = "test reference 2"

# Active cell moves from A4 to A3
$ws.Range("A3").Select()
